# "Working Edit Vendor details" -- mark the Vendor-details checklist item
# as done and clear out a stray leftover checkbox value in the totals row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 = "Allow agents to manage vendor details." -- tick it off (was FALSE).
$ws.Range("A12").Value = $true

# Row 23 is a blank spacer row in the totals block; it had a stray FALSE
# checkbox value left in A23 that should just be cleared out.
$ws.Range("A23").Value = ""

# C24's SUMIF total recalculates automatically from the A-column checkboxes.

# Reflect where the user was working: scrolled down and selected E20.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E20").Select() | Out-Null
